$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neonatal_Deaths_Data")

# Correct the recorded neonatal death count for 2016 (D8): 84367 -> 84368
$ws.Range("D8").Value = 84368

# Move/save the selection onto the corrected cell (D8), replacing the
# previous full-row selection on row 15
$ws.Activate()
$ws.Range("D8").Select()
